$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate() | Out-Null

# --- G10: "LM051更新" -> "L7205更新" -------------------------------------
$ws.Range("G10").Value = "L7205更新"

# --- G11: "LM052更新" -> "LM051更新" -------------------------------------
$ws.Range("G11").Value = "LM051更新"

# --- G12: "LM052更新" -> "L7205更新", and pick up the full-border style
#     that G10 already has (was missing its top border before).
$ws.Range("G12").Value = "L7205更新"
$ws.Range("G10").Copy() | Out-Null
$ws.Range("G12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Move the saved cursor/selection to H6 (was D15)
$ws.Range("H6").Select() | Out-Null
